$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.361.45"
$ws.Range("E2").Value = "  -2.83%  "

$ws.Range("D3").Value = "3.690.52"
$ws.Range("E3").Value = "  -3.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "686.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.22%  "

$ws.Range("D7").Value = "3.688.98"
$ws.Range("E7").Value = "  -3.29%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -6.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -8.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -10.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.03%  "

$ws.Range("D14").Value = "4.308.89"
$ws.Range("E14").Value = "  -3.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.52"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -11.34%  "

$ws.Range("D16").Value = "3.684.21"
$ws.Range("E16").Value = "  -3.05%  "

$ws.Range("D17").Value = "69.432.18"
$ws.Range("E17").Value = "  -2.78%  "

$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -9.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -10.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "474.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.647"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.08%  "

$ws.Range("D25").Value = "3.831.05"
$ws.Range("E25").Value = "  -3.24%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -11.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -13.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.20"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -11.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -10.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.76"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.67"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -9.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -10.70%  "

$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.70"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.159"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.20"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -11.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.10"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.27"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0904"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -10.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.942"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "165.41"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.82"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -15.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "28.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000274"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.86"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.28%  "
